$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 64

# Columns A-D hold text-like values (date/time/weekday/week stored as text
# in this sheet, not real Excel dates/numbers). Force text interpretation
# the way a user pre-formatting cells as Text would, then clear the
# formatting residue so the cell keeps the sheet's default style.
$textCols = 1,2,3,4
$textVals = "2023-06-20","21:38:44","Tuesday","25"
for ($i = 0; $i -lt $textCols.Length; $i++) {
    $cell = $ws.Cells.Item($row, $textCols[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $textVals[$i]
    $cell.ClearFormats()
}

$ws.Cells.Item($row, 5).Value = 122188
$ws.Cells.Item($row, 6).Value = 133651
$ws.Cells.Item($row, 7).Value = 162501
$ws.Cells.Item($row, 8).Value = 133256
$ws.Cells.Item($row, 9).Value = 177331
$ws.Cells.Item($row, 10).Value = 114645
$ws.Cells.Item($row, 11).Value = 201716
$ws.Cells.Item($row, 12).Value = 225493
$ws.Cells.Item($row, 13).Value = 175469
$ws.Cells.Item($row, 14).Value = 103872
$ws.Cells.Item($row, 15).Value = 39319
$ws.Cells.Item($row, 16).Value = 33895
$ws.Cells.Item($row, 17).Value = 51968
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36415
$ws.Cells.Item($row, 20).Value = -1
